# Actualización automática 2025-07-04 14:05:09
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M7").Value = 651.5

# --- Sheet: VENTA MENSUAL ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F7").Value = 2250.87
$wsMensual.Range("F22").Value = 10283.78

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D16").Value = 7329.64
$wsCumpl.Range("E16").Value = 36936.6
$wsCumpl.Range("F16").Value = 0.1655808128271116

$wsCumpl.Range("D19").Value = 10283.78
$wsCumpl.Range("E19").Value = 55094.21762291768
$wsCumpl.Range("F19").Value = 0.1572972616768414
